# Applies the "esquema_miniumax" fixture schedule fix:
# clears the stale match-schedule cells and writes the corrected
# team/opponent values (column A team order is unchanged).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fixture")

# Remove the old (now incorrect) fixture entries
$ws.Range("F2").ClearContents()
$ws.Range("G2").ClearContents()
$ws.Range("R2").ClearContents()
$ws.Range("S2").ClearContents()
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("L3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("J5").ClearContents()
$ws.Range("Q5").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("J6").ClearContents()
$ws.Range("N6").ClearContents()
$ws.Range("R6").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("J7").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("B9").ClearContents()
$ws.Range("C9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("M10").ClearContents()
$ws.Range("Q10").ClearContents()
$ws.Range("C11").ClearContents()
$ws.Range("G11").ClearContents()
$ws.Range("J11").ClearContents()
$ws.Range("R11").ClearContents()

# Write the corrected fixture entries
$ws.Range("B2").Value = "@CHI"
$ws.Range("K2").Value = "CHI"
$ws.Range("G3").Value = "URU"
$ws.Range("P3").Value = "@URU"
$ws.Range("D4").Value = "@VEN"
$ws.Range("M4").Value = "VEN"
$ws.Range("B5").Value = "ARG"
$ws.Range("E5").Value = "VEN"
$ws.Range("K5").Value = "@ARG"
$ws.Range("P5").Value = "@VEN"
$ws.Range("I9").Value = "@URU"
$ws.Range("R9").Value = "URU"
$ws.Range("G10").Value = "@BOL"
$ws.Range("I10").Value = "PER"
$ws.Range("P10").Value = "BOL"
$ws.Range("R10").Value = "@PER"
$ws.Range("D11").Value = "BRA"
$ws.Range("E11").Value = "@CHI"
$ws.Range("M11").Value = "@BRA"
$ws.Range("P11").Value = "CHI"

